# Update InsideBet Data: Automatizado
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 4 (Excelsior)
$ws.Range("D4").Value = 45.7
$ws.Range("E4").Value = 24
$ws.Range("F4").Value = 264
$ws.Range("G4").Value = 2160
$ws.Range("H4").Value = 24
$ws.Range("I4").Value = 26
$ws.Range("J4").Value = 16
$ws.Range("K4").Value = 42
$ws.Range("L4").Value = 24
$ws.Range("O4").Value = 36
$ws.Range("Q4").Value = 1.08
$ws.Range("R4").Value = 0.67
$ws.Range("S4").Value = 1.75
$ws.Range("U4").Value = 1.67

# Row 6 (Fortuna Sittard)
$ws.Range("D6").Value = 45.5
$ws.Range("E6").Value = 24
$ws.Range("F6").Value = 264
$ws.Range("G6").Value = 2160
$ws.Range("H6").Value = 24
$ws.Range("I6").Value = 35
$ws.Range("J6").Value = 26
$ws.Range("K6").Value = 61
$ws.Range("L6").Value = 32
$ws.Range("M6").Value = 3
$ws.Range("N6").Value = 3
$ws.Range("O6").Value = 53
$ws.Range("Q6").Value = 1.46
$ws.Range("R6").Value = 1.08
$ws.Range("S6").Value = 2.54
$ws.Range("T6").Value = 1.33
$ws.Range("U6").Value = 2.42

$wb.Save()
